$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Section 1: resistance test table (rows 18-33) ---
# Observed values entered for the K column; formulas in L/O/P recalc automatically.
$rows1 = 18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33
foreach ($r in $rows1) {
    $ws.Cells.Item($r, 11).Value = 0.7   # column K
}

# --- Section rows 40-47 (second resistance/voltage block) ---
# GND-referenced rows get 0.3, VA/VD rows get 0
$ws.Cells.Item(40, 11).Value = 0.3
$ws.Cells.Item(41, 11).Value = 0.3
$ws.Cells.Item(42, 11).Value = 0
$ws.Cells.Item(43, 11).Value = 0.3
$ws.Cells.Item(44, 11).Value = 0.3
$ws.Cells.Item(45, 11).Value = 0
$ws.Cells.Item(46, 11).Value = 0.3
$ws.Cells.Item(47, 11).Value = 0.3

# --- Section: acceptance checkboxes (M57:M59) ---
$ws.Range("M57").Value = "y"
$ws.Range("M58").Value = "y"
$ws.Range("M59").Value = "y"

# --- Section: overall result for section 2 (row 60) ---
$ws.Range("B60").Value = "pass"

# --- Section: HV / voltage drop measurements (rows 66-71) ---
$ws.Range("K66").Value = 98.8
$ws.Range("K67").Value = 9.877
$ws.Range("C70").Value = -0.0175
$ws.Range("C71").Value = -0.0173

# --- Final decision (row 80): tester name and test date ---
$ws.Range("B80").Value = "Amanda"
$testDate = Get-Date -Year 2014 -Month 10 -Day 6 -Hour 0 -Minute 0 -Second 0
$ws.Range("F80").Value = $testDate.Date
